$p = $ppt.ActivePresentation

# --- Slide 3 ("Project Process Flow") ---
$s3 = $p.Slides.Item(3)

# Un-hide the slide (removes show="0" from the slide XML)
$s3.SlideShowTransition.Hidden = 0

# "Built a website using GitHub Pages" -> "...and Hugo"
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shape = $s3.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.TextRange.Text -eq "Built a website using GitHub Pages") {
            $shape.TextFrame.TextRange.Text = "Built a website using GitHub Pages and Hugo"
        }
    }
}

# --- Slide 5 ("Understanding AI Image Generation") ---
$s5 = $p.Slides.Item(5)

# "Used Google’s Go language..." -> "Used Python and Google’s Go language..."
# (COM TextRange.Text normalizes the curly apostrophe to a straight one when read back,
# so search with a straight apostrophe but write the proper typographic one.)
$needle = "Used Google's Go language to build the background infrastructure for a potential blogging website"
$replacement = "Used Python and Google’s Go language to build the background infrastructure for a potential blogging website"
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $shape = $s5.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf($needle)
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, $needle.Length)
            $sub.Text = $replacement
        }
    }
}

# --- Slide 7 ("DEMO TIME") ---
$s7 = $p.Slides.Item(7)

# Un-hide the slide (removes show="0" from the slide XML)
$s7.SlideShowTransition.Hidden = 0
